$d = $word.ActiveDocument

# Build the WordprocessingML for the three new paragraphs that must be
# inserted immediately before the document's final (empty) paragraph,
# followed by that same empty paragraph (preserved, since InsertXML
# replaces the target range's contents).
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$rFonts = '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>'

$para1 = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:rPr>' + $rFonts + '</w:rPr></w:pPr>' + `
    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t xml:space="preserve">, highlighting the capacity of these methods to understand observed fire behavior from satellite-derived information harmonized with wall-to-wall forest </w:t></w:r>' + `
    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>inventory</w:t></w:r>' + `
    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t xml:space="preserve"> data</w:t></w:r>' + `
  '</w:p>'

$para2 = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:rPr>' + $rFonts + '</w:rPr></w:pPr>' + `
  '</w:p>'

$bodyText = 'The influence of species total live basal area had a much lower or insignificant effect on CBIbc for all forest types, including aspen. Average aspen tree height and diameter had diverging effects on both FRPc and CBIbc. In aspen forests, average tree height had a significant negative influence with a -1.9% to -2.3% average decrease in FRPc and CBIbc, respectively, for each unit increase in tree height. Conversely, tree diameter had less pronounced but positive effect, where greater average diameter tended to increase both FRPc and CBIbc. For both responses, the gridcell average canopy cover percent had a strong positive effect, total dead trees abundance had a weak negative effect, and the gridcell diversity of species contributing to live basal area had a significant positive effect ('

$para3 = '<w:p ' + $wNs + '>' + `
    '<w:pPr><w:spacing w:before="120" w:after="240"/><w:rPr>' + $rFonts + '</w:rPr></w:pPr>' + `
    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>' + $bodyText + '</w:t></w:r>' + `
    '<w:r><w:rPr>' + $rFonts + '<w:b/><w:bCs/></w:rPr><w:t>Figure SX, Table SX</w:t></w:r>' + `
    '<w:r><w:rPr>' + $rFonts + '</w:rPr><w:t>).</w:t></w:r>' + `
  '</w:p>'

$trailingEmptyPara = '<w:p ' + $wNs + '></w:p>'

$xml = $para1 + $para2 + $para3 + $trailingEmptyPara

# The final paragraph in the document body is the empty one right before
# the section properties. Replacing its range with the XML above keeps an
# equivalent trailing empty paragraph while inserting the three new ones
# ahead of it.
$target = $d.Paragraphs.Last
[void]$target.Range.InsertXML($xml)
